$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New cell L2: shared string ", 5,10, 50,100, 500"
$ws.Range("L2").Value = ", 5,10, 50,100, 500"

# New column M formulas: =6*K<row> for rows 9, 11, 12, 14
$ws.Range("M9").Formula = "=6*K9"
$ws.Range("M11").Formula = "=6*K11"
$ws.Range("M12").Formula = "=6*K12"
$ws.Range("M14").Formula = "=6*K14"

# Update the selection/view state to match the saved file
$ws.Range("O14").Select()
